$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.224.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2844"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06541"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.845.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6902"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.085"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "266.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.210.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.120.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.235"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.152"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.471"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.930"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09910"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.343"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.454"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04734"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.781"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.259"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.933"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4149"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8325"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "979.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.072"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.141"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05653"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
